$d = $word.ActiveDocument

$replacements = @(
    @{old = "155÷3="; new = "762÷7="},
    @{old = "345÷5="; new = "543÷4="},
    @{old = "473÷3="; new = "535÷7="},
    @{old = "365÷6="; new = "368÷9="},
    @{old = "972÷8="; new = "253÷5="},
    @{old = "253÷3="; new = "324÷2="},
    @{old = "489÷4="; new = "252÷4="},
    @{old = "653÷8="; new = "229÷3="},
    @{old = "896÷7="; new = "222÷8="},
    @{old = "186÷5="; new = "417÷9="},
    @{old = "692÷4="; new = "565÷2="},
    @{old = "416÷5="; new = "176÷2="},
    @{old = "420÷8="; new = "415÷8="},
    @{old = "278÷9="; new = "295÷6="},
    @{old = "884÷5="; new = "908÷5="},
    @{old = "978÷6="; new = "742÷3="},
    @{old = "178÷9="; new = "986÷7="},
    @{old = "750÷7="; new = "435÷4="},
    @{old = "351÷2="; new = "881÷7="},
    @{old = "868÷6="; new = "966÷6="},
    @{old = "800÷6="; new = "631÷3="},
    @{old = "354÷4="; new = "694÷3="},
    @{old = "581÷9="; new = "606÷3="},
    @{old = "269÷9="; new = "504÷9="},
    @{old = "462÷5="; new = "450÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
